$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Mean (H) and Std (I) values per row, as described in the commit diff
$updates = @(
    @(26, 0.71949, 0.02054),
    @(27, 0.40996, 0.05847),
    @(28, 0.72111, 0.02283),
    @(29, 0.42923, 0.07077),
    @(30, 0.71936, 0.02036),
    @(31, 0.40765, 0.05976),
    @(32, 0.72412, 0.0229),
    @(33, 0.4327, 0.07371999999999999),
    @(34, 0.6522, 0.02122),
    @(35, 0.17458, 0.04481),
    @(36, 0.65182, 0.02122),
    @(37, 0.17574, 0.04594),
    @(38, 0.6521400000000001, 0.02095),
    @(39, 0.17458, 0.04481),
    @(40, 0.65188, 0.02108),
    @(41, 0.17574, 0.04594),
    @(66, 0.65507, 0.02985),
    @(67, 0.24929, 0.05557),
    @(68, 0.67179, 0.02994),
    @(69, 0.33875, 0.05841),
    @(70, 0.65603, 0.02909),
    @(71, 0.24968, 0.05528),
    @(72, 0.6750699999999999, 0.03001),
    @(73, 0.34877, 0.05399),
    @(74, 0.63468, 0.02645),
    @(75, 0.16612, 0.03893),
    @(76, 0.6398200000000001, 0.02792),
    @(77, 0.19734, 0.04356),
    @(78, 0.635, 0.02635),
    @(79, 0.16612, 0.03893),
    @(80, 0.64059, 0.0284),
    @(81, 0.1985, 0.04438),
    @(106, 0.71243, 0.01965),
    @(107, 0.38849, 0.04555),
    @(108, 0.7148099999999999, 0.02061),
    @(109, 0.40774, 0.05256),
    @(110, 0.71275, 0.0198),
    @(111, 0.38656, 0.04855),
    @(112, 0.7160300000000001, 0.02062),
    @(113, 0.40851, 0.05232),
    @(114, 0.64822, 0.02302),
    @(115, 0.17416, 0.05232),
    @(116, 0.64854, 0.02281),
    @(117, 0.17493, 0.05237),
    @(118, 0.64822, 0.02302),
    @(119, 0.17416, 0.05232),
    @(120, 0.64854, 0.02281),
    @(121, 0.17493, 0.05237),
    @(146, 0.68886, 0.02645),
    @(147, 0.30906, 0.05841),
    @(148, 0.69803, 0.03009),
    @(149, 0.38649, 0.07124999999999999),
    @(150, 0.68809, 0.02713),
    @(151, 0.30906, 0.05841),
    @(152, 0.70394, 0.03084),
    @(153, 0.39689, 0.07421),
    @(154, 0.65716, 0.02777),
    @(155, 0.19158, 0.05444),
    @(156, 0.65787, 0.02749),
    @(157, 0.21123, 0.06074),
    @(158, 0.65729, 0.02815),
    @(159, 0.19158, 0.05444),
    @(160, 0.65832, 0.02759),
    @(161, 0.21084, 0.05995),
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 8).Value = $u[1]
    $ws.Cells.Item($r, 9).Value = $u[2]
}

$wb.Save()